# Updated cryptos list on Mon Oct  9 13:24:41 UTC 2023 with GitHub Actions
#
# This script mirrors the data-refresh diff: Price (column D) and
# Volume(1h) (column E) are updated for the rows whose figures changed,
# and rows 26/27 (Cosmos / BinanceUSD) swap ranking positions with
# refreshed values. All D/E cells are plain text in the workbook, so
# a leading apostrophe is used when assigning values that look
# numeric, to keep them stored as text (matching the original
# inlineStr cell type) rather than being auto-converted to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with new data ---
# A leading apostrophe forces Excel to treat the assigned value as literal text,
# matching the original inlineStr (text) cell type for columns D and E.

$ws.Range("D2").Value = "'27.527.40"
$ws.Range("E2").Value = "'  -0.98%  "

$ws.Range("D3").Value = "'1.594.64"
$ws.Range("E3").Value = "'  -1.95%  "

$ws.Range("E4").Value = "'  +0.34%  "

$ws.Range("D5").Value = "'207.78"
$ws.Range("E5").Value = "'  -1.35%  "

$ws.Range("E6").Value = "'  -3.61%  "

$ws.Range("E7").Value = "'  +0.34%  "

$ws.Range("D8").Value = "'22.25"
$ws.Range("E8").Value = "'  -4.21%  "

$ws.Range("D10").Value = "'0.0592"
$ws.Range("E10").Value = "'  -3.41%  "

$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "'  -0.91%  "

$ws.Range("D12").Value = "'1.820.37"
$ws.Range("E12").Value = "'  -1.85%  "

$ws.Range("D13").Value = "'1.605.40"
$ws.Range("E13").Value = "'  -1.28%  "

$ws.Range("E14").Value = "'  -3.90%  "

$ws.Range("E15").Value = "'  -4.24%  "

$ws.Range("D16").Value = "'63.36"
$ws.Range("E16").Value = "'  -3.03%  "

$ws.Range("D17").Value = "'27.538.03"
$ws.Range("E17").Value = "'  -0.89%  "

$ws.Range("D18").Value = "'218.60"
$ws.Range("E18").Value = "'  -4.63%  "

$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "'  -3.14%  "

$ws.Range("D20").Value = "'0.0₃0695"
$ws.Range("E20").Value = "'  -3.55%  "

$ws.Range("E22").Value = "'  -2.37%  "

$ws.Range("D23").Value = "'9.64"
$ws.Range("E23").Value = "'  -4.21%  "

$ws.Range("E24").Value = "'  -2.21%  "

$ws.Range("D25").Value = "'154.73"
$ws.Range("E25").Value = "'  +0.46%  "

# --- Rows 26 and 27 swapped places in the ranking (Cosmos moved above BinanceUSD) ---
# Row 26 becomes Cosmos (was BinanceUSD); Row 27 becomes BinanceUSD (was Cosmos).
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'6.75"
$ws.Range("E26").Value = "'  -2.02%  "

$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  +0.34%  "

$ws.Range("D28").Value = "'15.04"
$ws.Range("E28").Value = "'  -3.06%  "

$ws.Range("E29").Value = "'  -4.14%  "

$ws.Range("E30").Value = "'  -1.15%  "

$ws.Range("D31").Value = "'0.0469"
$ws.Range("E31").Value = "'  -2.55%  "

$ws.Range("D32").Value = "'3.27"
$ws.Range("E32").Value = "'  -4.13%  "

$ws.Range("D33").Value = "'1.361.09"
$ws.Range("E33").Value = "'  -2.28%  "

$ws.Range("E34").Value = "'  -4.36%  "

$ws.Range("E35").Value = "'  -2.39%  "

$ws.Range("D36").Value = "'0.962"
$ws.Range("E36").Value = "'  -4.88%  "

$ws.Range("E37").Value = "'  -0.97%  "

$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = "'  -2.38%  "

$ws.Range("D39").Value = "'0.539"
$ws.Range("E39").Value = "'  -2.53%  "

$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "'  -3.86%  "

$ws.Range("E41").Value = "'  +0.26%  "

$ws.Range("D42").Value = "'0.968"
$ws.Range("E42").Value = "'  -3.19%  "

$ws.Range("D43").Value = "'5.38"
$ws.Range("E43").Value = "'  -1.02%  "

$ws.Range("E44").Value = "'  -2.49%  "

$ws.Range("D45").Value = "'1.77"
$ws.Range("E45").Value = "'  -3.22%  "

$ws.Range("D46").Value = "'1.731.05"
$ws.Range("E46").Value = "'  -1.82%  "

$ws.Range("E47").Value = "'  -2.72%  "

$ws.Range("D48").Value = "'87.66"
$ws.Range("E48").Value = "'  -0.08%  "

$ws.Range("D49").Value = "'0.0₇0999"
$ws.Range("E49").Value = "'  +3.76%  "

$ws.Range("D50").Value = "'0.0969"
$ws.Range("E50").Value = "'  -3.97%  "

$ws.Range("E51").Value = "'  -0.98%  "
